# Update the NATMI edge-weight table with newly computed TPM-based values.
# Columns A-D (Sending cluster, Ligand symbol, Receptor symbol, Target cluster)
# and columns E, I, J, K, L are unchanged; only F, G, H, M, N, O, P, Q, R, S, T
# are refreshed per row (rows 2-7, one row per target cluster).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2 = @{ F = 0.5; G = 0.2098545; H = 0.419709;
           M = 20.79761; N = 41.59522;
           O = 0.1221346813081061; P = 0.08652867365760288;
           Q = 4.364472047745; R = 17.45788819098;
           S = 0.1221346813081061; T = 0.08652867365760288 }
    3 = @{ F = 0.5; G = 0.2098545; H = 0.419709;
           M = 13.56641266666666; N = 40.69923799999999;
           O = 0.0796692258166966; P = 0.08466480242237233;
           Q = 2.846972746956999; R = 17.081836481742;
           S = 0.0796692258166966; T = 0.08466480242237233 }
    4 = @{ F = 0.5; G = 0.2098545; H = 0.419709;
           M = 49.80229833333333; N = 149.406895;
           O = 0.2924657129041698; P = 0.3108044736787241;
           Q = 10.4512364155925; R = 62.707418493555;
           S = 0.2924657129041698; T = 0.3108044736787241 }
    5 = @{ F = 0.5; G = 0.2098545; H = 0.419709;
           M = 9.34483; N = 18.68966;
           O = 0.05487783615177078; P = 0.03887926283143963;
           Q = 1.961054627235; R = 7.84421850894;
           S = 0.05487783615177078; T = 0.03887926283143963 }
    6 = @{ F = 0.5; G = 0.2098545; H = 0.419709;
           M = 37.76134866666666; N = 113.284046;
           O = 0.2217548211149075; P = 0.2356597283761661;
           Q = 7.924388943768999; R = 47.546333662614;
           S = 0.2217548211149075; T = 0.2356597283761661 }
    7 = @{ F = 0.5; G = 0.2098545; H = 0.419709;
           M = 39.011729; N = 117.035187;
           O = 0.2290977227043493; P = 0.2434630590336949;
           Q = 8.1867868834305; R = 49.120721300583;
           S = 0.2290977227043493; T = 0.2434630590336949 }
}

foreach ($row in $updates.Keys) {
    $cols = $updates[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
}
